$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Init_amort ("model_Sim" data update, lag AV2016): multiply the Amount
#    columns (F, H, J) for rows 8-15 by 1000 (expressing AV2016lag data in
#    the same units as elsewhere in the model).
# ---------------------------------------------------------------------------
$wsAmort = $wb.Worksheets.Item("Init_amort")

$rows = 8..15
foreach ($r in $rows) {
    foreach ($col in @("F", "H", "J")) {
        $addr = "$col$r"
        $wsAmort.Range($addr).Value = $wsAmort.Range($addr).Value2 * 1000
    }
}

# Widen column J (10th column) a bit so the larger numbers remain readable.
$wsAmort.Columns.Item(10).ColumnWidth = 12.45

# ---------------------------------------------------------------------------
# 2. View-state bookkeeping to mirror the saved selections / zoom that were
#    captured for each sheet at the time of the edit.
# ---------------------------------------------------------------------------

# Init_amort: zoomed in further, cursor left on H16.
$wsAmort.Activate()
$excel.ActiveWindow.Zoom = 130
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 5
$wsAmort.Range("H16").Select()

# Init_unrecReturn: cursor left on D18.
$wsUnrec = $wb.Worksheets.Item("Init_unrecReturn")
$wsUnrec.Activate()
$wsUnrec.Range("D18").Select()

# Init_unrecReturn_raw: cursor left on G29.
$wsUnrecRaw = $wb.Worksheets.Item("Init_unrecReturn_raw")
$wsUnrecRaw.Activate()
$wsUnrecRaw.Range("G29").Select()

# Restore Init_amort as the active/selected sheet (matches tabSelected/activeTab
# in the saved workbook).
$wsAmort.Activate()
